$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cost")
$lo = $ws.ListObjects.Item(1)

# The query-table-backed Table_cost gains a new, unbound/calculated column.
# Adding via ListColumns first (Excel auto-names it), then set the header
# text explicitly so it lands on "Column1" and the table range grows to D63.
$newCol = $lo.ListColumns.Add()
$ws.Range("D1").Value = "Column1"

# New summary AVERAGE() formulas dropped into column G alongside several
# of the technology groups' last rows.
$ws.Range("G7").Formula  = "=AVERAGE(B3:B10)"
$ws.Range("G10").Formula = "=AVERAGE(B10:B22)"
$ws.Range("G23").Formula = "=AVERAGE(B23:B25)"
$ws.Range("G30").Formula = "=AVERAGE(B26:B30)"
$ws.Range("G33").Formula = "=AVERAGE(B31:B34,B23:B25)"
$ws.Range("G42").Formula = "=AVERAGE(B35:B55)"
$ws.Range("G55").Formula = "=AVERAGE(B56:B63)"

# The BECCS rows (56-63) lose their explicit per-cell style override,
# reverting tech/Source cells back to the default "Normal" style.
$ws.Range("A56:A63").Style = "Normal"
$ws.Range("C56:C63").Style = "Normal"

# Update the active selection/view to match the author's final state.
$ws.Activate() | Out-Null
$ws.Range("G30").Select() | Out-Null
